# Folder cleanup before last commit
# Add a "Total" row (row 14, under the task table) with sums for
# estimated time (D), ideal duration (F) and delta (G), label it in C14,
# and move the active selection to H15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Total label in C14: bold, italic, underlined, right-aligned ---
$cellC14 = $ws.Range("C14")
$cellC14.Value = "Total"
$cellC14.Font.Bold = $true
$cellC14.Font.Italic = $true
$cellC14.Font.Underline = $true
$cellC14.HorizontalAlignment = -4152   # xlRight
$cellC14.WrapText = $true

# --- Sum formulas in D14, F14, G14; bold + centered (E14 stays blank) ---
$ws.Range("D14").Formula = "=SUM(D2:D13)"
$ws.Range("F14").Formula = "=SUM(F2:F13)"
$ws.Range("G14").Formula = "=SUM(G2:G13)"

$totalsRange = $ws.Range("D14:G14")
$totalsRange.Font.Bold = $true
$totalsRange.HorizontalAlignment = -4108   # xlCenter
$totalsRange.WrapText = $true

# --- Narrow the shared "Delta" formula so it stops at row 13 (the last
#     data row before the new Total row) instead of reaching row 15 ---
$ws.Range("G3:G13").FormulaR1C1 = "=RC4-RC6"

# --- Move the active selection, matching the saved view state ---
$ws.Range("H15").Select()

$wb.Save()
